$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" values are textual numbers (e.g. "17.20", "2.560") that must
# keep their exact digits/trailing zeros. Force those cells to Text format
# before assigning so Excel does not silently renormalize them as numbers.

# Row 2
$ws.Range("D2").Value = "23.296.39"
$ws.Range("E2").Value = "  -0.22%  "

# Row 3
$ws.Range("D3").Value = "1.621.79"
$ws.Range("E3").Value = "  -0.22%  "

# Row 4
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.23%  "

# Row 5
$ws.Range("E5").Value = "  +0.31%  "

# Row 6
$ws.Range("D6").Value = "303.27"
$ws.Range("E6").Value = "  -0.34%  "

# Row 7
$ws.Range("D7").Value = "0.3789"
$ws.Range("E7").Value = "  +0.22%  "

# Row 8
$ws.Range("D8").Value = "51.98"
$ws.Range("E8").Value = "  +0.03%  "

# Row 9
$ws.Range("D9").Value = "0.3532"
$ws.Range("E9").Value = "  -2.42%  "

# Row 10
$ws.Range("D10").Value = "0.08076"
$ws.Range("E10").Value = "  -0.24%  "

# Row 11
$ws.Range("D11").Value = "1.206"
$ws.Range("E11").Value = "  -1.73%  "

# Row 13
$ws.Range("D13").Value = "21.98"
$ws.Range("E13").Value = "  -3.14%  "

# Row 14
$ws.Range("D14").Value = "6.372"
$ws.Range("E14").Value = "  -2.77%  "

# Row 15
$ws.Range("D15").Value = "7.196"
$ws.Range("E15").Value = "  -0.37%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001213"
$ws.Range("E16").Value = "  -2.84%  "

# Row 17
$ws.Range("D17").Value = "1.617.02"
$ws.Range("E17").Value = "  -0.45%  "

# Row 18
$ws.Range("D18").Value = "94.21"
$ws.Range("E18").Value = "  +0.76%  "

# Row 19
$ws.Range("D19").Value = "0.06927"
$ws.Range("E19").Value = "  +0.28%  "

# Row 20
$ws.Range("D20").Value = "6.511"
$ws.Range("E20").Value = "  +1.39%  "

# Row 21
$ws.Range("D21").Value = "1.002"
$ws.Range("E21").Value = "  +0.18%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.20"
$ws.Range("E22").Value = "  -3.77%  "

# Row 23
$ws.Range("D23").Value = "12.24"
$ws.Range("E23").Value = "  -3.54%  "

# Row 24
$ws.Range("D24").Value = "23.281.10"
$ws.Range("E24").Value = "  -0.26%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.560"
$ws.Range("E25").Value = "  +4.74%  "

# Row 26
$ws.Range("D26").Value = "3.085"
$ws.Range("E26").Value = "  -4.87%  "

# Row 27
$ws.Range("D27").Value = "20.85"
$ws.Range("E27").Value = "  -1.12%  "

# Row 28
$ws.Range("D28").Value = "151.25"
$ws.Range("E28").Value = "  +0.88%  "

# Row 29
$ws.Range("D29").Value = "5.254"
$ws.Range("E29").Value = "  -0.61%  "

# Row 30
$ws.Range("D30").Value = "131.87"
$ws.Range("E30").Value = "  -1.68%  "

# Row 31
$ws.Range("D31").Value = "1.798.24"
$ws.Range("E31").Value = "  -0.43%  "

# Row 32
$ws.Range("D32").Value = "1.067"
$ws.Range("E32").Value = "  +12.09%  "

# Row 33
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "6.449"
$ws.Range("E33").Value = "  -5.10%  "

# Row 34
$ws.Range("B34").Value = "WEMIXTOKEN"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.110"
$ws.Range("E34").Value = "  -8.44%  "

# Row 35
$ws.Range("D35").Value = "11.38"
$ws.Range("E35").Value = "  +3.40%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02710"
$ws.Range("E36").Value = "  -2.46%  "

# Row 37
$ws.Range("B37").Value = "Algorand"
$ws.Range("C37").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D37").Value = "0.2463"
$ws.Range("E37").Value = "  -2.03%  "

# Row 38
$ws.Range("B38").Value = "Stellar"
$ws.Range("C38").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D38").Value = "0.08677"
$ws.Range("E38").Value = "  -1.68%  "

# Row 39
$ws.Range("D39").Value = "0.06893"
$ws.Range("E39").Value = "  -3.28%  "

# Row 40
$ws.Range("D40").Value = "5.847"
$ws.Range("E40").Value = "  -3.93%  "

# Row 41
$ws.Range("D41").Value = "0.6874"

# Row 42
$ws.Range("D42").Value = "1.311"
$ws.Range("E42").Value = "  -3.67%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.90"
$ws.Range("E43").Value = "  -3.28%  "

# Row 44
$ws.Range("D44").Value = "15.26"
$ws.Range("E44").Value = "  -5.76%  "

# Row 45
$ws.Range("D45").Value = "1.002"
$ws.Range("E45").Value = "  +0.44%  "

# Row 46
$ws.Range("D46").Value = "0.6291"
$ws.Range("E46").Value = "  -2.47%  "

# Row 47
$ws.Range("E47").Value = "  -0.95%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.250"
$ws.Range("E48").Value = "  -2.86%  "

# Row 49
$ws.Range("D49").Value = "0.07881"
$ws.Range("E49").Value = "  -1.31%  "

# Row 50
$ws.Range("D50").Value = "127.86"
$ws.Range("E50").Value = "  +1.74%  "

# Row 51
$ws.Range("D51").Value = "1.171"
$ws.Range("E51").Value = "  -2.19%  "
